$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.628928422927856
$ws.Range("B1").Value = 2.598803043365479
$ws.Range("C1").Value = 3.034539222717285
$ws.Range("D1").Value = 3.045434236526489
$ws.Range("E1").Value = 1.06139075756073
